$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "An enterprise can be thought of as the overall business, made up of all the individual sites or workplaces."
$newText = "Represented here are enterprises, which can be thought of as the overall business, made up of all the individual sites or workplaces."

# Update the three cells in column C that reference the old text (rows 12-14)
$ws.Range("C12").Value = $newText
$ws.Range("C13").Value = $newText
$ws.Range("C14").Value = $newText

# Update the view/selection state to match the final saved state as closely
# as this host supports (scroll position + multi-cell selection anchored on
# the edited range).
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12:C14").Select()
